$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly price observation per row for "Haba" at the
# "Mercado Mayorista Lo Valledor de Santiago" market. This edit inserts one
# new, more-recent observation at row 230, pushing the existing rows
# 230..321 down to 231..322 (dimension grows from A1:R321 to A1:R322).
$ws.Rows(230).Insert()

$ws.Range("A230").Value = 6
$ws.Range("B230").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C230").Value = "Metropolitana"
$ws.Range("D230").Value = 44839
$ws.Range("E230").Value = 13
$ws.Range("F230").Value = 100112026
$ws.Range("G230").Value = "Haba"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 700
$ws.Range("K230").Value = 10000
$ws.Range("L230").Value = 12000
$ws.Range("M230").Value = 11343
$ws.Range("N230").Value = "`$/saco 25 kilos"
$ws.Range("O230").Value = "Región Metropolitana"
$ws.Range("P230").Value = 454
$ws.Range("Q230").Value = 25
$ws.Range("R230").Value = "Hortaliza"
